$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 9.361143666666667
$ws.Range("H2").Value = 28.083431
$ws.Range("I2").Value = 0.1965934861218525
$ws.Range("J2").Value = 0.1965934861218526
$ws.Range("M2").Value = 139.2986196666667
$ws.Range("N2").Value = 417.895859
$ws.Range("O2").Value = 0.6137320738580456
$ws.Range("P2").Value = 0.6137320738580456
$ws.Range("Q2").Value = 1303.994391268026
$ws.Range("R2").Value = 11735.94952141223
$ws.Range("S2").Value = 0.1206557279445475
$ws.Range("T2").Value = 0.1206557279445475
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 9.361143666666667
$ws.Range("H3").Value = 28.083431
$ws.Range("I3").Value = 0.1965934861218525
$ws.Range("J3").Value = 0.1965934861218526
$ws.Range("O3").Value = 0.09768689432339951
$ws.Range("P3").Value = 0.09768689432339951
$ws.Range("Q3").Value = 207.5550027837858
$ws.Range("R3").Value = 1867.995025054072
$ws.Range("S3").Value = 0.01920460710345412
$ws.Range("T3").Value = 0.01920460710345412
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 9.361143666666667
$ws.Range("H4").Value = 28.083431
$ws.Range("I4").Value = 0.1965934861218525
$ws.Range("J4").Value = 0.1965934861218526
$ws.Range("M4").Value = 35.78898466666666
$ws.Range("N4").Value = 107.366954
$ws.Range("O4").Value = 0.1576817331952585
$ws.Range("P4").Value = 0.1576817331952585
$ws.Range("Q4").Value = 335.0258271487971
$ws.Range("R4").Value = 3015.232444339174
$ws.Range("S4").Value = 0.03099920162659171
$ws.Range("T4").Value = 0.03099920162659172
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 9.361143666666667
$ws.Range("H5").Value = 28.083431
$ws.Range("I5").Value = 0.1965934861218525
$ws.Range("J5").Value = 0.1965934861218526
$ws.Range("M5").Value = 29.710182
$ws.Range("N5").Value = 89.13054600000001
$ws.Range("O5").Value = 0.1308992986232963
$ws.Range("P5").Value = 0.1308992986232963
$ws.Range("Q5").Value = 278.121282064814
$ws.Range("R5").Value = 2503.091538583326
$ws.Range("S5").Value = 0.02573394944725924
$ws.Range("T5").Value = 0.02573394944725924
$ws.Range("I6").Value = 0.2809659460057553
$ws.Range("J6").Value = 0.2809659460057554
$ws.Range("M6").Value = 139.2986196666667
$ws.Range("N6").Value = 417.895859
$ws.Range("O6").Value = 0.6137320738580456
$ws.Range("P6").Value = 0.6137320738580456
$ws.Range("Q6").Value = 1863.632539186631
$ws.Range("R6").Value = 16772.69285267968
$ws.Range("S6").Value = 0.1724378127255999
$ws.Range("T6").Value = 0.1724378127255999
$ws.Range("I7").Value = 0.2809659460057553
$ws.Range("J7").Value = 0.2809659460057554
$ws.Range("O7").Value = 0.09768689432339951
$ws.Range("P7").Value = 0.09768689432339951
$ws.Range("S7").Value = 0.02744669067593819
$ws.Range("T7").Value = 0.0274466906759382
$ws.Range("I8").Value = 0.2809659460057553
$ws.Range("J8").Value = 0.2809659460057554
$ws.Range("M8").Value = 35.78898466666666
$ws.Range("N8").Value = 107.366954
$ws.Range("O8").Value = 0.1576817331952585
$ws.Range("P8").Value = 0.1576817331952585
$ws.Range("Q8").Value = 478.8096000438096
$ws.Range("R8").Value = 4309.286400394286
$ws.Range("S8").Value = 0.04430319733503292
$ws.Range("T8").Value = 0.04430319733503293
$ws.Range("I9").Value = 0.2809659460057553
$ws.Range("J9").Value = 0.2809659460057554
$ws.Range("M9").Value = 29.710182
$ws.Range("N9").Value = 89.13054600000001
$ws.Range("O9").Value = 0.1308992986232963
$ws.Range("P9").Value = 0.1308992986232963
$ws.Range("Q9").Value = 397.483205884246
$ws.Range("R9").Value = 3577.348852958215
$ws.Range("S9").Value = 0.03677824526918431
$ws.Range("T9").Value = 0.03677824526918432
$ws.Range("G10").Value = 7.684952333333334
$ws.Range("H10").Value = 23.054857
$ws.Range("I10").Value = 0.1613917725961189
$ws.Range("J10").Value = 0.1613917725961189
$ws.Range("M10").Value = 139.2986196666667
$ws.Range("N10").Value = 417.895859
$ws.Range("O10").Value = 0.6137320738580456
$ws.Range("P10").Value = 0.6137320738580456
$ws.Range("Q10").Value = 1070.503252237463
$ws.Range("R10").Value = 9634.529270137164
$ws.Range("S10").Value = 0.09905130729904214
$ws.Range("T10").Value = 0.09905130729904214
$ws.Range("G11").Value = 7.684952333333334
$ws.Range("H11").Value = 23.054857
$ws.Range("I11").Value = 0.1613917725961189
$ws.Range("J11").Value = 0.1613917725961189
$ws.Range("O11").Value = 0.09768689432339951
$ws.Range("P11").Value = 0.09768689432339951
$ws.Range("Q11").Value = 170.3905377093982
$ws.Range("R11").Value = 1533.514839384584
$ws.Range("S11").Value = 0.01576586103426319
$ws.Range("T11").Value = 0.01576586103426319
$ws.Range("G12").Value = 7.684952333333334
$ws.Range("H12").Value = 23.054857
$ws.Range("I12").Value = 0.1613917725961189
$ws.Range("J12").Value = 0.1613917725961189
$ws.Range("M12").Value = 35.78898466666666
$ws.Range("N12").Value = 107.366954
$ws.Range("O12").Value = 0.1576817331952585
$ws.Range("P12").Value = 0.1576817331952585
$ws.Range("Q12").Value = 275.0366412217309
$ws.Range("R12").Value = 2475.329770995578
$ws.Range("S12").Value = 0.02544853442641105
$ws.Range("T12").Value = 0.02544853442641106
$ws.Range("G13").Value = 7.684952333333334
$ws.Range("H13").Value = 23.054857
$ws.Range("I13").Value = 0.1613917725961189
$ws.Range("J13").Value = 0.1613917725961189
$ws.Range("M13").Value = 29.710182
$ws.Range("N13").Value = 89.13054600000001
$ws.Range("O13").Value = 0.1308992986232963
$ws.Range("P13").Value = 0.1308992986232963
$ws.Range("Q13").Value = 228.3213324846581
$ws.Range("R13").Value = 2054.891992361922
$ws.Range("S13").Value = 0.0211260698364025
$ws.Range("T13").Value = 0.0211260698364025
$ws.Range("G14").Value = 17.19197166666666
$ws.Range("H14").Value = 51.57591499999999
$ws.Range("I14").Value = 0.3610487952762732
$ws.Range("J14").Value = 0.3610487952762733
$ws.Range("M14").Value = 139.2986196666667
$ws.Range("N14").Value = 417.895859
$ws.Range("O14").Value = 0.6137320738580456
$ws.Range("P14").Value = 0.6137320738580456
$ws.Range("Q14").Value = 2394.817922515109
$ws.Range("R14").Value = 21553.36130263598
$ws.Range("S14").Value = 0.2215872258888561
$ws.Range("T14").Value = 0.2215872258888561
$ws.Range("G15").Value = 17.19197166666666
$ws.Range("H15").Value = 51.57591499999999
$ws.Range("I15").Value = 0.3610487952762732
$ws.Range("J15").Value = 0.3610487952762733
$ws.Range("O15").Value = 0.09768689432339951
$ws.Range("P15").Value = 0.09768689432339951
$ws.Range("Q15").Value = 381.1798914954977
$ws.Range("R15").Value = 3430.61902345948
$ws.Range("S15").Value = 0.03526973550974401
$ws.Range("T15").Value = 0.03526973550974401
$ws.Range("G16").Value = 17.19197166666666
$ws.Range("H16").Value = 51.57591499999999
$ws.Range("I16").Value = 0.3610487952762732
$ws.Range("J16").Value = 0.3610487952762733
$ws.Range("M16").Value = 35.78898466666666
$ws.Range("N16").Value = 107.366954
$ws.Range("O16").Value = 0.1576817331952585
$ws.Range("P16").Value = 0.1576817331952585
$ws.Range("Q16").Value = 615.283210368101
$ws.Range("R16").Value = 5537.548893312909
$ws.Range("S16").Value = 0.05693079980722283
$ws.Range("T16").Value = 0.05693079980722285
$ws.Range("G17").Value = 17.19197166666666
$ws.Range("H17").Value = 51.57591499999999
$ws.Range("I17").Value = 0.3610487952762732
$ws.Range("J17").Value = 0.3610487952762733
$ws.Range("M17").Value = 29.710182
$ws.Range("N17").Value = 89.13054600000001
$ws.Range("O17").Value = 0.1308992986232963
$ws.Range("P17").Value = 0.1308992986232963
$ws.Range("Q17").Value = 510.77660715551
$ws.Range("R17").Value = 4596.98946439959
$ws.Range("S17").Value = 0.04726103407045026
$ws.Range("T17").Value = 0.04726103407045027
